$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 14).Value = -40929.669
$ws.Cells.Item(17, 12).Value = 40593.669
$ws.Cells.Item(17, 8).Value = 13531.223
$ws.Cells.Item(17, 10).Value = 13531.223

$ws.Cells.Item(40, 14).Value = -3823.75
$ws.Cells.Item(40, 9).Value = 2591.6667
$ws.Cells.Item(40, 11).Value = 2591.6667
$ws.Cells.Item(40, 12).Value = 3473.75
$ws.Cells.Item(40, 8).Value = 2944.5
$ws.Cells.Item(40, 10).Value = 3473.75
$ws.Cells.Item(40, 13).Value = -2416.6667

$ws.Cells.Item(62, 14).Value = -3181.3334
$ws.Cells.Item(62, 9).Value = 5536.875
$ws.Cells.Item(62, 11).Value = 5536.875
$ws.Cells.Item(62, 12).Value = 1933.3334
$ws.Cells.Item(62, 8).Value = 4554.091
$ws.Cells.Item(62, 10).Value = 1933.3334
$ws.Cells.Item(62, 13).Value = -4912.875

$ws.Cells.Item(65, 14).Value = -15906.667
$ws.Cells.Item(65, 9).Value = 5536.875
$ws.Cells.Item(65, 11).Value = 27684.375
$ws.Cells.Item(65, 12).Value = 9666.666999999999
$ws.Cells.Item(65, 8).Value = 4554.091
$ws.Cells.Item(65, 10).Value = 1933.3334
$ws.Cells.Item(65, 13).Value = -24564.375

$ws.Cells.Item(69, 14).Value = -13646
$ws.Cells.Item(69, 9).Value = 3000
$ws.Cells.Item(69, 11).Value = 9000
$ws.Cells.Item(69, 12).Value = 11898
$ws.Cells.Item(69, 8).Value = 3724.5
$ws.Cells.Item(69, 10).Value = 3966
$ws.Cells.Item(69, 13).Value = -8126

$ws.Cells.Item(72, 14).Value = -44430
$ws.Cells.Item(72, 9).Value = 3000
$ws.Cells.Item(72, 11).Value = 27000
$ws.Cells.Item(72, 12).Value = 35694
$ws.Cells.Item(72, 8).Value = 3724.5
$ws.Cells.Item(72, 10).Value = 3966
$ws.Cells.Item(72, 13).Value = -22632

$ws.Cells.Item(107, 14).Value = -4140
$ws.Cells.Item(107, 9).Value = 854.25
$ws.Cells.Item(107, 11).Value = 854.25
$ws.Cells.Item(107, 12).Value = 300
$ws.Cells.Item(107, 8).Value = 803.86365
$ws.Cells.Item(107, 10).Value = 300
$ws.Cells.Item(107, 13).Value = 1065.75

$ws.Cells.Item(129, 14).Value = -16658.7694
$ws.Cells.Item(129, 12).Value = 6658.769400000001
$ws.Cells.Item(129, 8).Value = 237371.39
$ws.Cells.Item(129, 10).Value = 2219.5898

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 14).Value = -100347.664
$ws.Cells.Item(6, 9).Value = 56667.668
$ws.Cells.Item(6, 11).Value = 56667.668
$ws.Cells.Item(6, 12).Value = 100001.664
$ws.Cells.Item(6, 8).Value = 71112.336
$ws.Cells.Item(6, 10).Value = 100001.664
$ws.Cells.Item(6, 13).Value = -56494.668

$ws.Cells.Item(61, 9).Value = 2153.4666
$ws.Cells.Item(61, 11).Value = 2153.4666
$ws.Cells.Item(61, 8).Value = 2288.3235
$ws.Cells.Item(61, 13).Value = -1941.4666

$ws.Cells.Item(74, 14).Value = -3965.9333
$ws.Cells.Item(74, 9).Value = 803.0714
$ws.Cells.Item(74, 11).Value = 803.0714
$ws.Cells.Item(74, 12).Value = 2217.9333
$ws.Cells.Item(74, 8).Value = 1534.8966
$ws.Cells.Item(74, 10).Value = 2217.9333
$ws.Cells.Item(74, 13).Value = 70.92859999999996

$ws.Cells.Item(77, 14).Value = -19825.6665
$ws.Cells.Item(77, 9).Value = 803.0714
$ws.Cells.Item(77, 11).Value = 4015.357
$ws.Cells.Item(77, 12).Value = 11089.6665
$ws.Cells.Item(77, 8).Value = 1534.8966
$ws.Cells.Item(77, 10).Value = 2217.9333
$ws.Cells.Item(77, 13).Value = 352.643

$ws.Cells.Item(97, 14).Value = -1578.6
$ws.Cells.Item(97, 9).Value = 554.4286
$ws.Cells.Item(97, 11).Value = 554.4286
$ws.Cells.Item(97, 12).Value = 586.6
$ws.Cells.Item(97, 8).Value = 559.30304
$ws.Cells.Item(97, 10).Value = 586.6
$ws.Cells.Item(97, 13).Value = -58.42859999999996

$ws.Cells.Item(122, 14).Value = -11664.4
$ws.Cells.Item(122, 9).Value = 6006
$ws.Cells.Item(122, 11).Value = 18018
$ws.Cells.Item(122, 12).Value = 6764.400000000001
$ws.Cells.Item(122, 8).Value = 3326.5715
$ws.Cells.Item(122, 10).Value = 2254.8
$ws.Cells.Item(122, 13).Value = -15568

$ws.Cells.Item(136, 9).Value = 2153.4666
$ws.Cells.Item(136, 11).Value = 6460.399800000001
$ws.Cells.Item(136, 8).Value = 2288.3235
$ws.Cells.Item(136, 13).Value = -3910.399800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(122, 14).Value = -50576
$ws.Cells.Item(122, 12).Value = 40776
$ws.Cells.Item(122, 8).Value = 40776
$ws.Cells.Item(122, 10).Value = 40776

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 14).ClearContents()

$ws.Cells.Item(58, 14).Value = -3384.6155
$ws.Cells.Item(58, 9).Value = 1632.2931
$ws.Cells.Item(58, 11).Value = 1632.2931
$ws.Cells.Item(58, 12).Value = 2978.6155
$ws.Cells.Item(58, 8).Value = 1878.8029
$ws.Cells.Item(58, 10).Value = 2978.6155
$ws.Cells.Item(58, 13).Value = -1429.2931

$ws.Cells.Item(136, 14).Value = -14035.8465
$ws.Cells.Item(136, 9).Value = 1632.2931
$ws.Cells.Item(136, 11).Value = 4896.879300000001
$ws.Cells.Item(136, 12).Value = 8935.8465
$ws.Cells.Item(136, 8).Value = 1878.8029
$ws.Cells.Item(136, 10).Value = 2978.6155
$ws.Cells.Item(136, 13).Value = -2346.879300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 9).Value = 500000060
$ws.Cells.Item(9, 11).Value = 1500000180
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 8).Value = 500000060
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 13).Value = -1499999956
$ws.Cells.Item(9, 14).ClearContents()

$ws.Cells.Item(38, 14).Value = -1001211.52
$ws.Cells.Item(38, 9).Value = 179.14285
$ws.Cells.Item(38, 11).Value = 537.4285500000001
$ws.Cells.Item(38, 12).Value = 1000517.52
$ws.Cells.Item(38, 8).Value = 154022.23
$ws.Cells.Item(38, 10).Value = 333505.84
$ws.Cells.Item(38, 13).Value = -190.4285500000001

$ws.Cells.Item(131, 14).Value = -78496.56299999999
$ws.Cells.Item(131, 9).Value = 14821.286
$ws.Cells.Item(131, 11).Value = 44463.858
$ws.Cells.Item(131, 12).Value = 68416.56299999999
$ws.Cells.Item(131, 8).Value = 22240.98
$ws.Cells.Item(131, 10).Value = 22805.521
$ws.Cells.Item(131, 13).Value = -39423.858

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 8).Value = 5550000
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 14).ClearContents()

$ws.Cells.Item(122, 14).Value = -9662.5
$ws.Cells.Item(122, 9).Value = 1941.625
$ws.Cells.Item(122, 11).Value = 5824.875
$ws.Cells.Item(122, 12).Value = 4762.5
$ws.Cells.Item(122, 8).Value = 1823.5834
$ws.Cells.Item(122, 10).Value = 1587.5
$ws.Cells.Item(122, 13).Value = -3374.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 14).Value = -3349
$ws.Cells.Item(7, 9).Value = 2778.1538
$ws.Cells.Item(7, 11).Value = 2778.1538
$ws.Cells.Item(7, 12).Value = 3125
$ws.Cells.Item(7, 8).Value = 2910.2856
$ws.Cells.Item(7, 10).Value = 3125
$ws.Cells.Item(7, 13).Value = -2666.1538

$ws.Cells.Item(100, 9).Value = 1879.8
$ws.Cells.Item(100, 11).Value = 1879.8
$ws.Cells.Item(100, 8).Value = 1999.8334
$ws.Cells.Item(100, 13).Value = -1338.8

$ws.Cells.Item(126, 14).Value = -14315
$ws.Cells.Item(126, 9).Value = 2778.1538
$ws.Cells.Item(126, 11).Value = 8334.4614
$ws.Cells.Item(126, 12).Value = 9375
$ws.Cells.Item(126, 8).Value = 2910.2856
$ws.Cells.Item(126, 10).Value = 3125
$ws.Cells.Item(126, 13).Value = -5864.4614

$ws.Cells.Item(133, 14).Value = -38139.625
$ws.Cells.Item(133, 12).Value = 33079.625
$ws.Cells.Item(133, 8).Value = 33079.625
$ws.Cells.Item(133, 10).Value = 33079.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 14).Value = -4201.2666
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 2953.2666
$ws.Cells.Item(62, 8).Value = 2953.2666
$ws.Cells.Item(62, 10).Value = 2953.2666
$ws.Cells.Item(62, 13).ClearContents()

$ws.Cells.Item(65, 14).Value = -21006.333
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 14766.333
$ws.Cells.Item(65, 8).Value = 2953.2666
$ws.Cells.Item(65, 10).Value = 2953.2666
$ws.Cells.Item(65, 13).ClearContents()

$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 14).ClearContents()

$ws.Cells.Item(113, 14).Value = -7538
$ws.Cells.Item(113, 9).Value = 1043.1111
$ws.Cells.Item(113, 11).Value = 3129.3333
$ws.Cells.Item(113, 12).Value = 3198
$ws.Cells.Item(113, 8).Value = 1048.8334
$ws.Cells.Item(113, 10).Value = 1066
$ws.Cells.Item(113, 13).Value = -959.3333000000002

$ws.Cells.Item(132, 14).Value = -18180.5
$ws.Cells.Item(132, 9).Value = 881.44116
$ws.Cells.Item(132, 11).Value = 2644.32348
$ws.Cells.Item(132, 12).Value = 13120.5
$ws.Cells.Item(132, 8).Value = 1249.0264
$ws.Cells.Item(132, 10).Value = 4373.5
$ws.Cells.Item(132, 13).Value = -114.32348

$ws.Cells.Item(136, 14).Value = -12723.0768
$ws.Cells.Item(136, 9).Value = 144013.42
$ws.Cells.Item(136, 11).Value = 432040.26
$ws.Cells.Item(136, 12).Value = 7623.0768
$ws.Cells.Item(136, 8).Value = 24069.436
$ws.Cells.Item(136, 10).Value = 2541.0256
$ws.Cells.Item(136, 13).Value = -429490.26
